$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.386
$ws.Range("D7").Value = -7.505
$ws.Range("B9").Value = 5.680999999999999
$ws.Range("D12").Value = -7.542
$ws.Range("B13").Value = 5.302999999999999
$ws.Range("D14").Value = -7.778
$ws.Range("E15").Value = 15.94
$ws.Range("B16").Value = 4.974
$ws.Range("B18").Value = 5.236
$ws.Range("D19").Value = -7.851999999999999
$ws.Range("B20").Value = 7.162999999999999
$ws.Range("B26").Value = 5.444000000000001
$ws.Range("D26").Value = -7.798
$ws.Range("B27").Value = 6.090000000000001
$ws.Range("D27").Value = -8.086000000000002
$ws.Range("E28").Value = 17.099
$ws.Range("B29").Value = 5.673
$ws.Range("D29").Value = -7.452000000000001
$ws.Range("E33").Value = 17.213
$ws.Range("B35").Value = 8.439
$ws.Range("E35").Value = 16.517
$ws.Range("B36").Value = 7.923
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.258
$ws.Range("E38").Value = 16.792
$ws.Range("E43").Value = 17.127
$ws.Range("E44").Value = 16.93
$ws.Range("B45").Value = 5.895999999999999
$ws.Range("E45").Value = 17.058
$ws.Range("D47").Value = -7.480999999999999
$ws.Range("E47").Value = 16.984
$ws.Range("D51").Value = -8.575000000000001
$ws.Range("E51").Value = 16.616
$ws.Range("D52").Value = -7.874000000000001
$ws.Range("E54").Value = 16.47
$ws.Range("B55").Value = 5.185
$ws.Range("D55").Value = -7.943000000000001
$ws.Range("B57").Value = 5.327999999999999
$ws.Range("E57").Value = 16.641
$ws.Range("E62").Value = 16.3
$ws.Range("E63").Value = 17.471
$ws.Range("E67").Value = 16.887
$ws.Range("B69").Value = 5.581
$ws.Range("D69").Value = -7.325999999999999
$ws.Range("D70").Value = -7.449
$ws.Range("E70").Value = 17.468
$ws.Range("B76").Value = 6.495
$ws.Range("D76").Value = -7.839
$ws.Range("B78").Value = 8.550999999999998
$ws.Range("D81").Value = -7.785000000000001
$ws.Range("E81").Value = 17.015
$ws.Range("B82").Value = 5.289
$ws.Range("B83").Value = 5.129
$ws.Range("D83").Value = -8.516999999999999
$ws.Range("E88").Value = 16.327
$ws.Range("B93").Value = 6.130999999999999
$ws.Range("D94").Value = -7.456
$ws.Range("E96").Value = 16.783
$ws.Range("B97").Value = 5.632
$ws.Range("E99").Value = 16.553
$ws.Range("D100").Value = -8.276999999999999
$ws.Range("D102").Value = -7.865
